# Update existing rows 2-4 with the new office-title data, and append a
# new row 5 for the newly created official title (per commit message:
# "create new official titles").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: que fan dou wu shi -> di he shi
$ws.Range("A2").Value = 803799
$ws.Range("B2").Value = "堤河使"
$ws.Range("D2").Value = "Commissioner of Dams and rivers"
$ws.Range("E2").Value = "di he shi"

# Row 3: dao tian dou wu shi -> shen wei zuo xiang dou zhi hui shi
$ws.Range("A3").Value = 803800
$ws.Range("B3").Value = "神衛左廂都指揮使"
$ws.Range("D3").Value = "Commander-in-chief of the Left Inspired Guard Wing"
$ws.Range("E3").Value = "shen wei zuo xiang dou zhi hui shi"

# Row 4: que yan zhi zhi pan guan -> shen wei you xiang dou zhi hui shi
$ws.Range("A4").Value = 803801
$ws.Range("B4").Value = "神衛右廂都指揮使"
$ws.Range("D4").Value = "Commander-in-chief of the Right Inspired Guard Wing"
$ws.Range("E4").Value = "shen wei you xiang dou zhi hui shi"

# Row 5 (new): shi wei qin jun ma jun si dou zhi hui shi
$ws.Range("A5").Value = 803802
$ws.Range("B5").Value = "侍衛親軍馬軍司都指揮使"
# c_dy / c_source columns hold numeric-looking codes but are stored as text
# in this sheet (same as every other data row), so force Text format before
# assigning the value.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "15"
$ws.Range("D5").Value = "Commander-in-chief of the Metropolitan Cavalry"
$ws.Range("E5").Value = "shi wei qin jun ma jun si dou zhi hui shi"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "64847"
